$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.097.19'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '3.725.81'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.42'
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.32'
$ws.Range('E6').Value = '  +0.94%  '
$ws.Range('D7').Value = '3.723.81'
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.516'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  +0.89%  '
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('E13').Value = '  -0.76%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.11'
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('D15').Value = '4.357.96'
$ws.Range('D16').Value = '3.735.25'
$ws.Range('E16').Value = '  -0.28%  '
$ws.Range('D17').Value = '68.168.89'
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.85'
$ws.Range('E18').Value = '  -2.65%  '
$ws.Range('E19').Value = '  +0.26%  '
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('E21').Value = '  +2.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '464.94'
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000147'
$ws.Range('E24').Value = '  +10.98%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.69'
$ws.Range('E25').Value = '  +1.34%  '
$ws.Range('E26').Value = '  +0.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.83'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.03'
$ws.Range('E28').Value = '  -0.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').Value = '3.880.90'
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.75'
$ws.Range('E31').Value = '  -3.76%  '
$ws.Range('E32').Value = '  -0.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.65'
$ws.Range('E33').Value = '  +0.46%  '
$ws.Range('E34').Value = '  -0.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.16'
$ws.Range('E35').Value = '  +2.51%  '
$ws.Range('D37').Value = '3.684.20'
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('E38').Value = '  -0.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.40'
$ws.Range('E39').Value = '  -0.68%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.138'
$ws.Range('E40').Value = '  +1.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.995'
$ws.Range('E41').Value = '  +0.63%  '
$ws.Range('E42').Value = '  +1.20%  '
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '43.71'
$ws.Range('E45').Value = '  +15.23%  '
$ws.Range('E46').Value = '  -1.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.48'
$ws.Range('E47').Value = '  +2.71%  '
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.43'
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '144.28'
$ws.Range('E50').Value = '  -0.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '387.98'
$ws.Range('E51').Value = '  -0.13%  '
